$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9472426772117615
$ws.Range("B1").Value = 3.25377345085144
$ws.Range("C1").Value = 4.142512798309326
$ws.Range("D1").Value = 3.056813478469849
$ws.Range("E1").Value = 1.356643676757812
